$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "5/4/2015"
$ws.Range("D4").Value = "5/4/2015"
$ws.Range("D5").Value = "5/5/2015"
$ws.Range("D6").Value = "Fim de Checklist"
$ws.Range("D7").Value = "5/11/2015"
$ws.Range("D8").Value = "5/5/2015"

$ws.Range("D5").Select() | Out-Null
